# Update "想去人数" (want-to-go count) figures in column F
# for the "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

# Sheet "展览": row -> new F value
$wsExpo = $wb.Worksheets.Item("展览")
$expoUpdates = @{
    6  = 2283
    8  = 1712
    9  = 2890
    11 = 4301
    15 = 547
    16 = 255
    21 = 295
    22 = 4141
    24 = 3626
    25 = 1127
    27 = 535
    28 = 4365
    29 = 83
    30 = 481
    31 = 513
    32 = 455
}
foreach ($row in $expoUpdates.Keys) {
    $wsExpo.Range("F$row").Value = $expoUpdates[$row]
}

# Sheet "全部类型": row -> new F value
$wsAll = $wb.Worksheets.Item("全部类型")
$allUpdates = @{
    8  = 2283
    10 = 1712
    12 = 2890
    13 = 163
    14 = 4301
    18 = 547
    19 = 255
    25 = 295
    26 = 4141
    28 = 3626
    29 = 1127
    31 = 535
    32 = 4365
    33 = 83
    34 = 481
    35 = 513
    36 = 455
}
foreach ($row in $allUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allUpdates[$row]
}

$wb.Save()
